$d = $word.ActiveDocument

# wdHeaderFooterPrimary = 1
$section = $d.Sections(1)
$header = $section.Headers(1)
$header.Range.InsertAfter("Questionnaire 11")
$header.Range.Paragraphs(1).Style = "Header"
$header.Range.Paragraphs(1).Alignment = 1  # wdAlignParagraphCenter

$textRange = $header.Range.Paragraphs(1).Range
$textRange.MoveEnd(1, -1) | Out-Null  # exclude the paragraph mark
$textRange.Font.Name = "Arial"
$textRange.Font.Size = 12
